$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 420-421, shifting the existing rows 420:448 down to 422:450.
$ws.Range('A420:T421').Insert()

# New row 420 data
$ws.Range('A420').Value = 10
$ws.Range('B420').Value = 'Vega Modelo de Temuco'
$ws.Range('C420').Value = 'La Araucanía'
$ws.Range('D420').Value = 44585
$ws.Range('E420').Value = 9
$ws.Range('F420').Value = 'Fruta'
$ws.Range('G420').Value = 100108
$ws.Range('H420').Value = 'Tropicales y subtropicales'
$ws.Range('I420').Value = 100108006
$ws.Range('J420').Value = 'Plátano'
$ws.Range('K420').Value = 'Barraganete'
$ws.Range('L420').Value = 'Primera'
$ws.Range('M420').Value = 200
$ws.Range('N420').Value = 28000
$ws.Range('O420').Value = 28000
$ws.Range('P420').Value = 28000
$ws.Range('Q420').Value = '$/caja 20 kilos'
$ws.Range('R420').Value = 'Ecuador'
$ws.Range('S420').Value = 1400
$ws.Range('T420').Value = 20

# New row 421 data
$ws.Range('A421').Value = 10
$ws.Range('B421').Value = 'Vega Modelo de Temuco'
$ws.Range('C421').Value = 'La Araucanía'
$ws.Range('D421').Value = 44585
$ws.Range('E421').Value = 9
$ws.Range('F421').Value = 'Fruta'
$ws.Range('G421').Value = 100108
$ws.Range('H421').Value = 'Tropicales y subtropicales'
$ws.Range('I421').Value = 100108006
$ws.Range('J421').Value = 'Plátano'
$ws.Range('K421').Value = 'Sin especificar'
$ws.Range('L421').Value = 'Pintón'
$ws.Range('M421').Value = 1800
$ws.Range('N421').Value = 15000
$ws.Range('O421').Value = 16000
$ws.Range('P421').Value = 15444
$ws.Range('Q421').Value = '$/caja 20 kilos'
$ws.Range('R421').Value = 'Ecuador'
$ws.Range('S421').Value = 772
$ws.Range('T421').Value = 20
